$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New AgTests / AgPosit figures added for existing row 178 ---
$ws.Range("F178").Value = 46
$ws.Range("G178").Value = 0

# --- Revised AgTests / AgPosit figures for existing rows ---
$updates = @(
    @{ Row = 273; F = 31179; G = 1663 },
    @{ Row = 280; F = 34818; G = 2306 },
    @{ Row = 281; F = 45813; G = 3108 },
    @{ Row = 282; F = 44064; G = 2560 },
    @{ Row = 285; F = 40532; G = 3172 },
    @{ Row = 286; F = 52920; G = 4138 },
    @{ Row = 287; F = 56605; G = 3551 },
    @{ Row = 288; F = 57106; G = 3844 },
    @{ Row = 289; F = 59351; G = 3371 },
    @{ Row = 292; F = 80767; G = 7090 },
    @{ Row = 293; F = 80182; G = 5528 },
    @{ Row = 294; F = 89656; G = 4697 },
    @{ Row = 299; F = 64039; G = 6553 },
    @{ Row = 300; F = 70634; G = 6702 },
    @{ Row = 301; F = 70706; G = 5600 },
    @{ Row = 302; F = 77333; G = 5511 },
    @{ Row = 351; F = 150261 },
    @{ Row = 384; F = 172549 },
    @{ Row = 398; F = 300853 },
    @{ Row = 465; F = 61582; G = 56 },
    @{ Row = 479; F = 42609 },
    @{ Row = 512; F = 8613 },
    @{ Row = 518; F = 7189 },
    @{ Row = 519; F = 7988 },
    @{ Row = 520; F = 10333 },
    @{ Row = 521; F = 6822 },
    @{ Row = 522; F = 5149 },
    @{ Row = 523; F = 10204 },
    @{ Row = 524; F = 7818 },
    @{ Row = 525; F = 7610 },
    @{ Row = 526; F = 8775 },
    @{ Row = 527; F = 11450 },
    @{ Row = 528; F = 8020 },
    @{ Row = 529; F = 5626 },
    @{ Row = 530; F = 12657 },
    @{ Row = 531; F = 9174 },
    @{ Row = 532; F = 10187 },
    @{ Row = 533; F = 11790 },
    @{ Row = 534; F = 16636 },
    @{ Row = 535; F = 10050 },
    @{ Row = 536; F = 7891 },
    @{ Row = 537; F = 13516; G = 51 },
    @{ Row = 538; F = 11130 },
    @{ Row = 539; F = 10417; G = 47 },
    @{ Row = 540; F = 12317; G = 63 },
    @{ Row = 541; F = 16336; G = 66 },
    @{ Row = 542; F = 10172; G = 50 },
    @{ Row = 543; F = 4563; G = 31 }
)

foreach ($u in $updates) {
    $ws.Range("F" + $u.Row).Value = $u.F
    if ($u.ContainsKey("G")) {
        $ws.Range("G" + $u.Row).Value = $u.G
    }
}

# --- Brand new rows appended at the bottom (544-550) ---
$newRows = @(
    @{ Row = 544; A = 44438; B = 394923; C = 6779;  D = 132; E = 12548; F = 14156; G = 95 },
    @{ Row = 545; A = 44439; B = 395122; C = 6546;  D = 199; E = 12548; F = 16383; G = 109 },
    @{ Row = 546; A = 44440; B = 395300; C = 3468;  D = 178; E = 12548; F = 3668;  G = 51 },
    @{ Row = 547; A = 44441; B = 395532; C = 7555;  D = 232; E = 12549; F = 13552; G = 144 },
    @{ Row = 548; A = 44442; B = 395861; C = 11605; D = 329; E = 12549; F = 14820; G = 140 },
    @{ Row = 549; A = 44443; B = 396080; C = 21034; D = 219; E = 12549; F = 8212;  G = 55 },
    @{ Row = 550; A = 44444; B = 396181; C = 13726; D = 101; E = 12551; F = 5045;  G = 50 }
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
}
